# Apply the weekly cryptos list refresh (values scraped on Mon Sep  2 22:13:02 UTC 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume columns hold text (not numbers) in this sheet -- force Text format
# before writing so numeric-looking strings (e.g. "526.86") are not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.811.96"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "2.546.53"
$ws.Range("E3").Value = "  +3.26%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "526.86"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "134.41"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").Value = "2.546.12"
$ws.Range("E9").Value = "  +2.84%  "
$ws.Range("D10").Value = "0.0990"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("D12").Value = "5.20"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "2.996.23"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").Value = "58.828.32"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").Value = "22.44"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("E17").Value = "  +1.66%  "
$ws.Range("D18").Value = "2.547.13"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").Value = "10.76"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").Value = "324.66"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "4.21"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  +7.96%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "65.11"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "0.412"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "7.48"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("D29").Value = "0.0₃0759"
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("D30").Value = "1.74"
$ws.Range("E30").Value = "  +3.02%  "
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("D32").Value = "168.70"
$ws.Range("E32").Value = "  -1.32%  "
$ws.Range("D33").Value = "6.39"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "18.36"
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("D38").Value = "3.99"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  +3.61%  "
$ws.Range("D40").Value = "36.80"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").Value = "0.788"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").Value = "281.44"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "134.31"
$ws.Range("E44").Value = "  +8.90%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "5.09"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").Value = "0.605"
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("D47").Value = "0.0924"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").Value = "0.0507"
$ws.Range("E48").Value = "  +3.74%  "
$ws.Range("D49").Value = "17.91"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").Value = "0.0217"
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").Value = "17.19"
$ws.Range("E51").Value = "  +1.65%  "
